# Updated symbol list with refreshed coinranking.com market data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'306.43"
$ws.Range("E2").Value = "'-0.28%"

# Row 3
$ws.Range("D3").Value = "'38.94"
$ws.Range("E3").Value = "'6.94%"

# Row 4
$ws.Range("D4").Value = "'5.099"
$ws.Range("E4").Value = "'0.87%"

# Row 5
$ws.Range("D5").Value = "'0.08051"
$ws.Range("E5").Value = "'-0.63%"

# Row 6
$ws.Range("D6").Value = "'1.939"
$ws.Range("E6").Value = "'-10.09%"

# Row 7
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = "'8.001"
$ws.Range("E7").Value = "'2.05%"

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.9313"
$ws.Range("E8").Value = "'0.39%"

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = "'0.1457"
$ws.Range("E9").Value = "'1.36%"

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1926"
$ws.Range("E10").Value = "'-0.16%"

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.08993"
$ws.Range("E11").Value = "'-1.17%"

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.03502"
$ws.Range("E12").Value = "'1.41%"

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09780"
$ws.Range("E13").Value = "'-1.30%"

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001395"
$ws.Range("E14").Value = "'-1.02%"

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = "'0.005918"
$ws.Range("E15").Value = "'-6.23%"

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = "'3.785"
$ws.Range("E16").Value = "'-1.48%"

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = "'4.200"
$ws.Range("E17").Value = "'1.27%"

# Row 19
$ws.Range("D19").Value = "'0.3420"
$ws.Range("E19").Value = "'-0.98%"

# Row 20
$ws.Range("D20").Value = "'0.1303"
$ws.Range("E20").Value = "'0.09%"

# Row 21
$ws.Range("D21").Value = "'4.779"
$ws.Range("E21").Value = "'-1.14%"

# Row 22
$ws.Range("D22").Value = "'0.2415"
$ws.Range("E22").Value = "'3.29%"

# Row 23
$ws.Range("D23").Value = "'0.04387"
$ws.Range("E23").Value = "'0.60%"

# Row 25
$ws.Range("D25").Value = "'0.004279"
$ws.Range("E25").Value = "'-13.02%"

# Row 26
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'0.19%"

# Row 39
$ws.Range("D39").Value = "'0.02040"
$ws.Range("E39").Value = "'1.08%"

# Row 40
$ws.Range("D40").Value = "'0.05029"
$ws.Range("E40").Value = "'-2.73%"

# Row 41
$ws.Range("D41").Value = "'0.007436"
$ws.Range("E41").Value = "'-1.15%"

# Row 42
$ws.Range("D42").Value = "'0.01001"
$ws.Range("E42").Value = "'-1.40%"

# Row 43
$ws.Range("E43").Value = "'-1.20%"

# Row 44
$ws.Range("D44").Value = "'0.002122"
$ws.Range("E44").Value = "'-1.20%"

# Row 45
$ws.Range("D45").Value = "'0.009034"
$ws.Range("E45").Value = "'-9.29%"

# Row 46
$ws.Range("D46").Value = "'0.00006174"
$ws.Range("E46").Value = "'-1.65%"

# Row 47
$ws.Range("E47").Value = "'0.20%"

# Row 48
$ws.Range("D48").Value = "'0.002784"

# Row 49
$ws.Range("E49").Value = "'28.26%"

# Row 50
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.20%"

# Row 51
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.20%"

